# Updated schema test list
# Adds "correct" verdicts (and matching cell formatting) for a few rows of
# the test table on the first worksheet, and updates the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 19: E19/F19 get "correct" (format like E17/F17), G19 gets the same
#     fill but stays empty (format like G17) ---
$ws.Range("E17").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$ws.Range("E19").Value = "correct"

$ws.Range("F17").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$ws.Range("F19").Value = "correct"

$ws.Range("G17").Copy() | Out-Null
$ws.Range("G19").PasteSpecial(-4122) | Out-Null

# --- Row 21: E21/F21 get "correct" (format like E17/F17) ---
$ws.Range("E17").Copy() | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").Value = "correct"

$ws.Range("F17").Copy() | Out-Null
$ws.Range("F21").PasteSpecial(-4122) | Out-Null
$ws.Range("F21").Value = "correct"

# --- Row 22: E22 gets "correct" (format like E18), F22 gets " correct"
#     (leading space, default/no special formatting) ---
$ws.Range("E18").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "correct"

$ws.Range("F22").Value = " correct"

# --- Row 23: D23 keeps its value but switches to the white-fill format
#     already used by E27/F27 (style 24); E23/F23 get "correct"
#     (format like E17/F17), G23 gets the matching empty fill (like G17) ---
$ws.Range("E27").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null

$ws.Range("E17").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = "correct"

$ws.Range("F17").Copy() | Out-Null
$ws.Range("F23").PasteSpecial(-4122) | Out-Null
$ws.Range("F23").Value = "correct"

$ws.Range("G17").Copy() | Out-Null
$ws.Range("G23").PasteSpecial(-4122) | Out-Null

# --- Row 24: E24/F24 previously-empty cells now get "correct" using the
#     formatting already used by C20/D20 (style 21) ---
$ws.Range("C20").Copy() | Out-Null
$ws.Range("E24").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").Value = "correct"

$ws.Range("D20").Copy() | Out-Null
$ws.Range("F24").PasteSpecial(-4122) | Out-Null
$ws.Range("F24").Value = "correct"

# Clear clipboard marquee and restore the saved cursor position/selection.
$excel.CutCopyMode = $false
$ws.Range("E22").Select() | Out-Null
